$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's pair of observations (dates 44615, "Primera"/"Segunda" quality rows)
# is inserted at row 320, pushing all subsequent rows down by two.
$ws.Range("A320:A321").EntireRow.Insert()

# Row 320 - Calidad "Primera"
$ws.Cells.Item(320, 1).Value = 8
$ws.Cells.Item(320, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(320, 3).Value = "Coquimbo"
$ws.Cells.Item(320, 4).Value = 44615
$ws.Cells.Item(320, 4).NumberFormat = $ws.Cells.Item(322, 4).NumberFormat()
$ws.Cells.Item(320, 5).Value = 4
$ws.Cells.Item(320, 6).Value = 100112009
$ws.Cells.Item(320, 7).Value = "Acelga"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 2400
$ws.Cells.Item(320, 11).Value = 500
$ws.Cells.Item(320, 12).Value = 600
$ws.Cells.Item(320, 13).Value = 550
$ws.Cells.Item(320, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(320, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(320, 16).Value = 275
$ws.Cells.Item(320, 17).Value = 2
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# Row 321 - Calidad "Segunda"
$ws.Cells.Item(321, 1).Value = 8
$ws.Cells.Item(321, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(321, 3).Value = "Coquimbo"
$ws.Cells.Item(321, 4).Value = 44615
$ws.Cells.Item(321, 4).NumberFormat = $ws.Cells.Item(322, 4).NumberFormat()
$ws.Cells.Item(321, 5).Value = 4
$ws.Cells.Item(321, 6).Value = 100112009
$ws.Cells.Item(321, 7).Value = "Acelga"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Segunda"
$ws.Cells.Item(321, 10).Value = 1280
$ws.Cells.Item(321, 11).Value = 400
$ws.Cells.Item(321, 12).Value = 450
$ws.Cells.Item(321, 13).Value = 425
$ws.Cells.Item(321, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(321, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(321, 16).Value = 212
$ws.Cells.Item(321, 17).Value = 2
$ws.Cells.Item(321, 18).Value = "Hortaliza"

Write-Output "Done"
